{"js": "// Replace each two-digit-division expression in the worksheet table with\n// its new value, preserving run formatting (font, size, etc.).\n// Old -> [new values in document order] -- a couple of old values repeat\n// (\"28\u00f78=\" and \"71\u00f79=\") and map to two different replacements depending\n// on which occurrence (first vs. second) is being edited.\nconst replacements = [\n  { old: \"23\u00f77=\", news: [\"63\u00f73=\"] },\n  { old: \"47\u00f72=\", news: [\"77\u00f79=\"] },\n  { old: \"56\u00f78=\", news: [\"66\u00f78=\"] },\n  { old: \"66\u00f79=\", news: [\"38\u00f73=\"] },\n  { old: \"80\u00f72=\", news: [\"21\u00f77=\"] },\n  { old: \"21\u00f74=\", news: [\"85\u00f78=\"] },\n  { old: \"44\u00f75=\", news: [\"52\u00f75=\"] },\n  { old: \"28\u00f78=\", news: [\"83\u00f77=\", \"77\u00f75=\"] },\n  { old: \"71\u00f79=\", news: [\"49\u00f76=\", \"81\u00f78=\"] },\n  { old: \"65\u00f79=\", news: [\"50\u00f76=\"] },\n  { old: \"27\u00f75=\", news: [\"27\u00f73=\"] },\n  { old: \"59\u00f72=\", news: [\"96\u00f79=\"] },\n  { old: \"71\u00f78=\", news: [\"69\u00f79=\"] },\n  { old: \"88\u00f76=\", news: [\"22\u00f74=\"] },\n  { old: \"17\u00f79=\", news: [\"97\u00f72=\"] },\n  { old: \"97\u00f76=\", news: [\"15\u00f75=\"] },\n  { old: \"71\u00f72=\", news: [\"15\u00f72=\"] },\n  { old: \"78\u00f77=\", news: [\"83\u00f79=\"] },\n  { old: \"67\u00f75=\", news: [\"53\u00f79=\"] },\n  { old: \"32\u00f78=\", news: [\"26\u00f74=\"] },\n  { old: \"67\u00f73=\", news: [\"82\u00f74=\"] },\n  { old: \"50\u00f75=\", news: [\"94\u00f74=\"] },\n  { old: \"76\u00f74=\", news: [\"43\u00f72=\"] },\n];\n\nconst body = context.document.body;\n\n// Run one search per distinct old value, then walk the matches in\n// document order, pairing them up with the replacement list.\nconst searches = replacements.map(r =>\n  body.search(r.old, { matchCase: true, matchWholeWord: false })\n);\nsearches.forEach(s => s.load(\"items\"));\nawait context.sync();\n\nreplacements.forEach((r, i) => {\n  const items = searches[i].items;\n  if (items.length !== r.news.length) {\n    throw new Error(\n      `Expected ${r.news.length} match(es) for \"${r.old}\" but found ${items.length}`\n    );\n  }\n  items.forEach((range, j) => {\n    range.insertText(r.news[j], Word.InsertLocation.replace);\n  });\n});\n\nawait context.sync();\n", "ps1": "# Replace each two-digit-division expression in the worksheet table with its\n# new value, preserving run formatting (font, size, etc.). A couple of old\n# values repeat (\"28\u00f78=\" and \"71\u00f79=\") and map to two different replacements\n# depending on occurrence order, so matches are walked strictly top-to-bottom:\n# a single Range is re-used and its Find cursor only ever moves forward from\n# the end of the previous match.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    , @(\"23\u00f77=\", \"63\u00f73=\")\n    , @(\"47\u00f72=\", \"77\u00f79=\")\n    , @(\"56\u00f78=\", \"66\u00f78=\")\n    , @(\"66\u00f79=\", \"38\u00f73=\")\n    , @(\"80\u00f72=\", \"21\u00f77=\")\n    , @(\"21\u00f74=\", \"85\u00f78=\")\n    , @(\"44\u00f75=\", \"52\u00f75=\")\n    , @(\"28\u00f78=\", \"83\u00f77=\")\n    , @(\"71\u00f79=\", \"49\u00f76=\")\n    , @(\"65\u00f79=\", \"50\u00f76=\")\n    , @(\"27\u00f75=\", \"27\u00f73=\")\n    , @(\"59\u00f72=\", \"96\u00f79=\")\n    , @(\"71\u00f78=\", \"69\u00f79=\")\n    , @(\"88\u00f76=\", \"22\u00f74=\")\n    , @(\"17\u00f79=\", \"97\u00f72=\")\n    , @(\"97\u00f76=\", \"15\u00f75=\")\n    , @(\"28\u00f78=\", \"77\u00f75=\")\n    , @(\"71\u00f72=\", \"15\u00f72=\")\n    , @(\"78\u00f77=\", \"83\u00f79=\")\n    , @(\"67\u00f75=\", \"53\u00f79=\")\n    , @(\"32\u00f78=\", \"26\u00f74=\")\n    , @(\"67\u00f73=\", \"82\u00f74=\")\n    , @(\"50\u00f75=\", \"94\u00f74=\")\n    , @(\"71\u00f79=\", \"81\u00f78=\")\n    , @(\"76\u00f74=\", \"43\u00f72=\")\n)\n\n$rng = $d.Content\n$isFirst = $true\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    if ($isFirst) {\n        $isFirst = $false\n    } else {\n        $rng.Collapse([Microsoft.Office.Interop.Word.WdCollapseDirection]::wdCollapseEnd)\n        $rng.End = $d.Content.End\n    }\n\n    $rng.Find.ClearFormatting()\n    $rng.Find.Text = $oldText\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = [Microsoft.Office.Interop.Word.WdFindWrap]::wdFindStop\n    $found = $rng.Find.Execute()\n    if (-not $found) {\n        throw (\"Could not find occurrence of '{0}'\" -f $oldText)\n    }\n\n    $rng.Text = $newText\n}\n"}
